$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C2:C5) from 2023-11-03 to 2023-11-13
$newDate = Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
